# Apply the "card replacement / removal of card slots / fixed effects" edit
# to the card_stats workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Header: "HP/Defence" -> "Defence/HP"
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Defence/HP"

# ---------------------------------------------------------------------------
# 2. Card table rows (3-10): Attack (C) becomes numeric 0 where it used to be
#    the placeholder text "Nan" / "5 to 20", and HP/Defence (D) becomes a
#    "Defence/HP" formatted "x/y" text value instead of the old free text.
# ---------------------------------------------------------------------------

# Attack column fixes (replace placeholder text with numeric 0)
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0

# Give C6 a 2-decimal numeric format (it used to hold the text "5 to 20")
$ws.Range("C6").NumberFormat = "0.00"

# Defence/HP column: force text format ("@") and set the "x/y" values
$defenceHpRange = $ws.Range("D3:D10")
$defenceHpRange.NumberFormat = "@"

# Set D9 ("5/0") before the "0/0" cells so the shared-string table ends up
# with the same ordering as the saved workbook.
$ws.Range("D9").Value = "5/0"
$ws.Range("D3").Value = "0/0"
$ws.Range("D4").Value = "0/0"
$ws.Range("D5").Value = "0/0"
$ws.Range("D6").Value = "0/0"
$ws.Range("D7").Value = "0/0"
$ws.Range("D8").Value = "0/0"
$ws.Range("D10").Value = "0/0"

# ---------------------------------------------------------------------------
# 3. Selection / view bookkeeping to mirror the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("E14:E15").Select()
